$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ions")

# Insert two new columns (I:J) before the existing "conc_units" column,
# pushing it to K and making room for the new "Dp" / "Dp_units" columns.
$ws.Columns("I:J").Insert()

# Header row
$ws.Range("I1").Value = "Dp"
$ws.Range("J1").Value = "Dp_units"

# Dp values (particle diameter) per ion row, with matching units column.
$ws.Range("I2").Value = 1
$ws.Range("I3:I6").Value = 0.000002
$ws.Range("I3:I6").NumberFormat = $ws.Range("E3").NumberFormat

$ws.Range("J2:J6").Value = "cm^2/s"
